# This script applies the shared-strings content updates described by the
# commit: newly-resolved Elastic-search-backed references fill in previously
# "Unknown Title/Abstract/[]  /not found/N/A/1970-01-01" placeholder rows, and
# a handful of author-list strings get re-written with corrected whitespace.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2
$ws.Range("E2").Value = '[Clara%González-Sanguino%NULL%0,                          Berta%Ausín%NULL%1,                          Miguel Ángel%Castellanos%NULL%1,                          Jesús%Saiz%NULL%1,                          Aída%López-Gómez%NULL%1,                          Carolina%Ugidos%NULL%1,                          Manuel%Muñoz%NULL%1]'

# E3
$ws.Range("E3").Value = '[Cuiyan%Wang%NULL%0,                          Riyu%Pan%NULL%0,                          Riyu%Pan%NULL%0,                          Xiaoyang%Wan%NULL%0,                          Yilin%Tan%NULL%0,                          Linkang%Xu%NULL%0,                          Cyrus S.%Ho%NULL%0,                          Roger C.%Ho%NULL%0,                          Roger C.%Ho%NULL%0]'

# E4
$ws.Range("E4").Value = '[Huiyao%Wang%NULL%1,                          Qian%Xia%NULL%1,                          Zhenzhen%Xiong%NULL%1,                          Zhixiong%Li%NULL%1,                          Weiyi%Xiang%NULL%1,                          Yiwen%Yuan%NULL%1,                          Yaya%Liu%NULL%1,                          Zhe%Li%NULL%1,                          Kenji%Hashimoto%NULL%0,                          Kenji%Hashimoto%NULL%0,                          Kenji%Hashimoto%NULL%0]'

# E5
$ws.Range("E5").Value = '[Yenan%Wang%xref no email%0,       Yu%Di%xref no email%0,       Junjie%Ye%xref no email%0,       Wenbin%Wei%xref no email%0]'

# E6
$ws.Range("E6").Value = '[Yuan%Wu%NULL%1,                          Jun%Wang%NULL%1,                          Chenggang%Luo%NULL%1,                          Sheng%Hu%NULL%1,                          Xi%Lin%NULL%1,                          Aimee E.%Anderson%NULL%1,                          Eduardo%Bruera%NULL%1,                          Xiaoxin%Yang%NULL%1,                          Shaozhong%Wei%NULL%1,                          Yu%Qian%173653835@qq.com%1]'

# E7
$ws.Range("E7").Value = '[Yingfei%Zhang%NULL%0,                          Zheng Feei%Ma%NULL%0,                          Zheng Feei%Ma%NULL%0]'

# C8
$ws.Range("C8").Value = '[Psychological status and sleep quality of nursing interns during the outbreak of COVID-19].'

# D8
$ws.Range("D8").Value = 'OBJECTIVE To investigate the psychological status and sleep quality of nursing interns in collective isolation during the outbreak of coronavirus disease 2019 (COVID-19) and provide evidence for adequate interventions.
 METHODS We surveyed a total of 95 nursing interns who were isolated collectively in a general teaching hospital in Guangzhou using a self-designed questionnaire, which consisted of a basic information form, self-rating anxiety scale (SAS), self-rating depression scale (SDS) and Pittsburgh Sleep Quality Index (PSQI).
 Descriptive analysis, single factor analysis and correlation analysis were used to analyze the current status of the interns'' psychology and sleep quality, the potential factors affecting their psychology and sleep quality, and the correlation between their psychological status and sleep quality.
 RESULTS The surveyed interns had SAS, SDS and PSQI score of 37.79±6.59, 43.98±9.74 and 5.20±3.14, respectively, which were significant higher than the national norms in China (P < 0.05).
 Correlation analysis indicated that both anxiety and depression were positively correlated with the sleep quality score (r=0.508 and 0.546, respectively).
 Univariate analysis showed that the major factors affecting the psychological status and sleep quality of the interns during collective isolation included recent contact with persons from the affected area before isolation and the onset of fever during the isolation.
 CONCLUSIONS These interns showed relatively high levels of anxiety and depression during the collective isolation to affect their sleep quality, and interventions should be timely administered to improve their mental health and sleep quality.
'

# E8
$ws.Range("E8").Value = '[ Xiaoyan%Sheng%elasticNoEmail%0,   Fenyu% Liu%elasticNoEmail%0,   Jin% Zhou%elasticNoEmail%0,   Rongrong% Liao%elasticNoEmail%0]'

# F8
$ws.Range("F8").Value = '10.12122/j.issn.1673-4254.2020.03.09'

# G8
$ws.Range("G8").Value = 'ELASTIC'

# H8
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = '2020-03-30'

# E9
$ws.Range("E9").Value = '[Cuiyan%Wang%NULL%0,                          Riyu%Pan%NULL%0,                          Xiaoyang%Wan%NULL%0,                          Yilin%Tan%NULL%0,                          Linkang%Xu%NULL%0,                          Roger S.%McIntyre%NULL%1,                          Faith N.%Choo%NULL%1,                          Bach%Tran%NULL%1,                          Roger%Ho%NULL%1,                          Vijay K.%Sharma%NULL%0,                          Cyrus%Ho%NULL%1]'

# E10
$ws.Range("E10").Value = '[Dana%Tzur Bitan%NULL%2,                          Ariella%Grossman-Giron%NULL%1,                          Yuval%Bloch%NULL%1,                          Yael%Mayer%NULL%1,                          Noga%Shiffman%NULL%1,                          Shlomo%Mendlovic%NULL%1]'

# C11
$ws.Range("C11").Value = '[Mental health survey of 230 medical staff in a tertiary infectious disease hospital for COVID-19].'

# D11
$ws.Range("D11").Value = 'Objective: To investigate the mental health of clinical first-line medical staff in COVID-19 epidemic and provide theoretical basis for psychological intervention.
 Method: The mental health status of the first-line medical staff was investigated by Self-rating Anxiety Acale (SAS) and Post-Traumatic Stress Disorder Self-rating Scale(PTSD-SS).
 From February 7 to 14, 2020, 246 medical staff were investigated who participated in the treatment of COVID-19 using cluster sampling , and received 230 responses, with a recovery rate of 93.5%.
 Results: The incidence of anxiety in medical staff was 23.04% (53/230), and the score of SAS was (42.91 ± 10.89).
 Among them, the incidence of severe anxiety, moderate anxiety and mild anxiety were 2.17% (5/230), 4.78% (11/230) and 16.09% (37/230), respectively.
 The incidence of anxiety in female medical staff was higher than that in male [25.67% (48/187) vs 11.63% (5/43), Z=-2.008, P=0.045], the score of SAS in female medical staff was higher than that in male [(43.78±11.12) vs (39.14 ± 9.01), t =-2.548, P=0.012].
 The incidence of anxiety in nurses was higher than that in doctors [26.88% (43/160) vs 14.29% (10/70), Z=-2.066, P=0.039], and the score of SAS in nurses was higher than that in doctors [(44.84±10.42) vs (38.50±10.72), t =-4.207, P<0.001].
 The incidence of stress disorder in medical staff was 27.39% (63/230), and the score of PTSD-SS was (42.92 ± 17.88).
 The score of PTSD-SS in female medical staff was higher than that of male [(44.30±18.42) vs(36.91 ± 13.95), t=-2.472, P=0.014].
 Conclusions: In COVID-19 epidemic, the incidence of anxiety and stress disorder is high among medical staff.
 Medical institutions should strengthen the training of psychological skills of medical staff.
 Special attention should be paid to the mental health of female nurses.
'

# E11
$ws.Range("E11").Value = '[ J Z%Huang%elasticNoEmail%0,   M F% Han%elasticNoEmail%0,   T D% Luo%elasticNoEmail%0,   A K% Ren%elasticNoEmail%0,   X P% Zhou%elasticNoEmail%0]'

# F11
$ws.Range("F11").Value = '10.3760/cma.j.cn121094-20200219-00063'

# G11
$ws.Range("G11").Value = 'ELASTIC'

# H11
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = '2020-03-04'

# E12
$ws.Range("E12").Value = '[Alison M.%Bacon%ambacon@plymouth.ac.uk%1,                          Philip J.%Corr%NULL%2,                          Philip J.%Corr%NULL%0]'

# E13
$ws.Range("E13").Value = '[Shuang-Jiang%Zhou%NULL%4,                          Li-Gang%Zhang%NULL%1,                          Lei-Lei%Wang%NULL%1,                          Zhao-Chang%Guo%NULL%1,                          Jing-Qi%Wang%NULL%1,                          Jin-Cheng%Chen%NULL%1,                          Mei%Liu%NULL%1,                          Xi%Chen%NULL%1,                          Jing-Xu%Chen%chenjx1110@163.com%1]'

# E14
$ws.Range("E14").Value = '[Haozheng%Cai%NULL%0,                          Baoren%Tu%NULL%0,                          Jing%Ma%NULL%0,                          Limin%Chen%NULL%0,                          Lei%Fu%NULL%0,                          Yongfang%Jiang%NULL%0,                          Quan%Zhuang%NULL%0]'

# E15
$ws.Range("E15").Value = '[Naiara%Ozamiz-Etxebarria%xref no email%1,       Maria%Dosil-Santamaria%xref no email%1,       Maitane%Picaza-Gorrochategui%xref no email%1,       Nahia%Idoiaga-Mondragon%xref no email%1]'

# E16
$ws.Range("E16").Value = '[Lei%Lei%NULL%1,                          Xiaoming%Huang%NULL%1,                          Shuai%Zhang%NULL%1,                          Jinrong%Yang%NULL%1,                          Lin%Yang%NULL%0,                          Min%Xu%NULL%1]'

# E17
$ws.Range("E17").Value = '[Shuai%Yuan%NULL%0,                          Zhenxin%Liao%NULL%0,                          Haojie%Huang%NULL%0,                          Boyue%Jiang%NULL%0,                          Xueyan%Zhang%NULL%0,                          Yingwen%Wang%NULL%0,                          Mingyi%Zhao%NULL%0]'

# E18
$ws.Range("E18").Value = '[Lijun%Kang%NULL%0,                          Simeng%Ma%NULL%0,                          Min%Chen%NULL%0,                          Jun%Yang%NULL%0,                          Ying%Wang%NULL%0,                          Ruiting%Li%NULL%0,                          Lihua%Yao%NULL%0,                          Hanping%Bai%NULL%0,                          Zhongxiang%Cai%NULL%0,                          Bing%Xiang Yang%NULL%0,                          Shaohua%Hu%NULL%0,                          Kerang%Zhang%NULL%0,                          Gaohua%Wang%NULL%0,                          Ci%Ma%NULL%0,                          Zhongchun%Liu%NULL%0]'

# E19
$ws.Range("E19").Value = '[Wenjun%Cao%wjcao16@czmc.edu.cn%0,                          Ziwei%Fang%NULL%0,                          Guoqiang%Hou%NULL%0,                          Mei%Han%NULL%0,                          Xinrong%Xu%NULL%0,                          Jiaxin%Dong%NULL%0,                          Jianzhong%Zheng%zhengjj@czmc.edu.cn%0]'

# E20
$ws.Range("E20").Value = '[Jianbo%Lai%NULL%0,                          Simeng%Ma%NULL%0,                          Ying%Wang%NULL%0,                          Zhongxiang%Cai%NULL%0,                          Jianbo%Hu%NULL%0,                          Ning%Wei%NULL%0,                          Jiang%Wu%NULL%0,                          Hui%Du%NULL%0,                          Tingting%Chen%NULL%0,                          Ruiting%Li%NULL%0,                          Huawei%Tan%NULL%0,                          Lijun%Kang%NULL%0,                          Lihua%Yao%NULL%0,                          Manli%Huang%NULL%0,                          Huafen%Wang%NULL%0,                          Gaohua%Wang%NULL%0,                          Zhongchun%Liu%NULL%0,                          Shaohua%Hu%NULL%0]'

# E21
$ws.Range("E21").Value = '[Yuanyuan%Mo%NULL%0,                          Lan%Deng%NULL%0,                          Lan%Deng%NULL%0,                          Liyan%Zhang%NULL%0,                          Qiuyan%Lang%NULL%0,                          Chunyan%Liao%NULL%0,                          Nannan%Wang%NULL%0,                          Mingqin%Qin%2026142822@qq.com%0,                          Huiqiao%Huang%820325832@qq.com%0]'

# E22
$ws.Range("E22").Value = '[Sijia%Li%NULL%0,                          Yilin%Wang%NULL%1,                          Jia%Xue%NULL%1,                          Nan%Zhao%NULL%1,                          Tingshao%Zhu%NULL%1]'

# E23
$ws.Range("E23").Value = '[Chen-Yun%Liu%NULL%0,                          Yun-zhi%Yang%NULL%1,                          Xiao-Ming%Zhang%NULL%1,                          Xinying%Xu%NULL%2,                          Xinying%Xu%NULL%0,                          Qing-Li%Dou%NULL%1,                          Wen-Wu%Zhang%NULL%1,                          Andy S. K.%Cheng%NULL%1]'

# E24
$ws.Range("E24").Value = '[Han%Xiao%NULL%0,                          Yan%Zhang%NULL%0,                          Desheng%Kong%NULL%0,                          Shiyue%Li%NULL%0,                          Ningxi%Yang%NULL%0]'

# E25
$ws.Range("E25").Value = '[Derek Gard-Ching%Ying%xref no email%1,       Sunfang%Jiang%xref no email%1,       Hua%Yang%xref no email%1,       Shanzhu%Zhu%xref no email%1]'

# E26
$ws.Range("E26").Value = '[Paolo%Soraci%paolo.soraci85@gmail.com%2,                          Ambra%Ferrari%a.ferrari50@campus.unimib.it%1,                          Francesco A.%Abbiati%francesco.abbiati@edu.unito.it%1,                          Elena%Del Fante%edelfante@hotmail.it%1,                          Rosanna%De Pace%Rosannadepace@libero.it%1,                          Antonino%Urso%antonino.urso@tin.it%1,                          Mark D.%Griffiths%mark.griffiths@ntu.ac.uk%1]'
